$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a new data row at position 12 (old row 12 "totals" becomes row 13,
#    old row 13 "footer" becomes row 14). Build its look from row 11 (the row
#    right above it) so fonts / borders / fills / merges all match the other
#    item rows in the table.
# ---------------------------------------------------------------------------
$ws.Rows("12:12").Insert()
$ws.Rows("12:12").RowHeight = 25.5

$ws.Range("A11:Q11").Copy()
$ws.Range("A12:Q12").PasteSpecial(-4122)

$ws.Range("A12:B12").Merge()
$ws.Range("C12:G12").Merge()
$ws.Range("H12:K12").Merge()
$ws.Range("L12:M12").Merge()
$ws.Range("N12:O12").Merge()

# Row 13 is now the totals row; it changes height from 25.5 to 24.75.
$ws.Rows("13:13").RowHeight = 24.75

# ---------------------------------------------------------------------------
# 2) Row 8 used to be "اولويز ماكس طويل جدا" -- replace it with the new
#    low-stock item that was added to the report.
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = "STRINGAZOLE 40MG 21 ENTERIC COATED TABLETS"
$ws.Range("H8").Value = "2:0"
$ws.Range("L8").Value = "'1"
$ws.Range("N8").Value = "126.00"
$ws.Range("P8").Value = "126.0000"

$ws.Range("L7").Copy()
$ws.Range("L8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Row 11 used to be "معطر جو FRIDA" -- it now holds a different item.
# ---------------------------------------------------------------------------
$ws.Range("C11").Value = "صوفي طويل جدا"
$ws.Range("H11").Value = "11:0"
$ws.Range("N11").Value = "50.00"
$ws.Range("P11").Value = "50.0000"

# ---------------------------------------------------------------------------
# 4) Row 12 (brand new) is item #6 -- "معطر جو FRIDA" reinstated with
#    different figures.
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "معطر جو FRIDA "
$ws.Range("H12").Value = "7:0"
$ws.Range("L12").Value = "'0"
$ws.Range("N12").Value = "65.00"
$ws.Range("P12").Value = "65.0000"

$ws.Range("L11").Copy()
$ws.Range("L12").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5) Totals row (now row 13): recompute sum of the price column (P7:P12).
# ---------------------------------------------------------------------------
$ws.Range("P13").Value = 367.2

# ---------------------------------------------------------------------------
# 6) Footer row (now row 14): refresh the generated-at timestamp.
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Sunday, 21 September, 2025 11:30 AM"
